$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 28, shifting existing rows 28-92 down to 29-93.
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with the new weekly record.
# Columns A, B, C, E, F, G, H, I, J, R keep the same values as the
# (now shifted) row below (old row 28), only the price / date / origin
# related columns change for this new entry.
$ws.Range("A28").Value = 4
$ws.Range("B28").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C28").Value = "Los Lagos"
$ws.Range("D28").Value = 44414
$ws.Range("E28").Value = 10
$ws.Range("F28").Value = 100112009
$ws.Range("G28").Value = "Acelga"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 200
$ws.Range("K28").Value = 4000
$ws.Range("L28").Value = 4000
$ws.Range("M28").Value = 4000
$ws.Range("N28").Value = "$/docena de atados (4 kilos)"
$ws.Range("O28").Value = "Región del Maule"
$ws.Range("P28").Value = 1000
$ws.Range("Q28").Value = 4
$ws.Range("R28").Value = "Hortaliza"
